$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# -----------------------------------------------------------------
# Change 1: delete the whole paragraph
#   "-      Ajouter des personnes avec le rôle de secrétaire"
# -----------------------------------------------------------------
$secMarker = "Ajouter des personnes avec le"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains($secMarker)) {
        $p.Range.Delete()
        break
    }
}

# -----------------------------------------------------------------
# Change 2: "Complete" + "r les pages<nbsp>: " -> "Compléter" + " les pages<nbsp>: "
# (keep the "calendrier," / " planning des " / "absences, convocations" / " "
#  runs that follow completely untouched)
# -----------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Completer les pages")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $start = $target.Range.Start

    # Rebuild the paragraph's leading runs from scratch by inserting each new
    # run at the (fixed/untouched) paragraph start -- done in reverse order so
    # the final reading order comes out correctly. Each InsertBefore() call on
    # a still-pristine paragraph start creates a brand-new run without
    # disturbing/merging the runs that already exist after it.
    $d.Range($start, $start).InsertBefore(" ")
    $d.Range($start, $start).InsertBefore("absences, convocations")
    $d.Range($start, $start).InsertBefore(" planning des ")
    $d.Range($start, $start).InsertBefore("calendrier,")
    $d.Range($start, $start).InsertBefore(" les pages" + $nbsp + ": ")
    $d.Range($start, $start).InsertBefore("Compléter")

    # Now remove the original (old) run content that follows what we just
    # inserted, all the way up to (but excluding) the paragraph mark.
    $newLen = 9 + 13 + 11 + 14 + 22 + 1
    $newEnd = $target.Range.End
    $oldStart = $start + $newLen
    $d.Range($oldStart, $newEnd - 1).Delete()
}

# -----------------------------------------------------------------
# Change 3: drop the parenthetical from the "Améliorer le visuel du site" bullet
# -----------------------------------------------------------------
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute(
    "Améliorer le visuel du site (dégager le stade notamment)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Améliorer le visuel du site ", 2) | Out-Null
